$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C (rows 2-357) from 45171 to 45172 (increment date by 1 day)
$range = $ws.Range("C2:C357")
foreach ($cell in $range.Cells) {
    if ($cell.Value2 -eq 45171) {
        $cell.Value2 = 45172
    }
}
